# STS IR Bot Performer / Data / Config.xlsx
# Commit: "Performer Process Ongoing. Review Sheet logic was added"
#
# Adds new Settings to the "Constants" sheet that support the new Review
# Sheet logic (template path, worksheet name, start cell, required columns,
# date format) plus a new "ReturnNameExceptionalCase" setting, and tidies
# up the stray fill style that used to live on B18.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Settings sheet: row-height touch-ups that come along with the edit
# (autofit height of the wrapped description cells).
# ---------------------------------------------------------------------
$wsSettings = $wb.Worksheets.Item("Settings")
$wsSettings.Rows.Item(3).RowHeight = 45
$wsSettings.Rows.Item(5).RowHeight = 30

# ---------------------------------------------------------------------
# Constants sheet: the bulk of the change.
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Constants")

# Same autofit row-height touch-ups as on the Settings sheet.
$ws.Rows.Item(2).RowHeight = 30
$ws.Rows.Item(3).RowHeight = 45

# B18 used to carry a leftover "applyFill" style (s="5"); drop it back to
# the default/Normal style so the cell is plain again.
$ws.Range("B18").Style = "Normal"

# Make room for the new Review Sheet settings: two brand new rows land
# right before the old "TaxFiledReport_SheetName" row (old row 26), which
# pushes everything from there on down by two rows.
$ws.Rows.Item(26).Resize(2).Insert()
$ws.Rows.Item(26).RowHeight = 14.25
$ws.Rows.Item(27).RowHeight = 14.25

# Fill in the new values. The order below matches the order the values
# were originally typed in (it controls shared-string allocation order),
# row 21 first, then row 19, then rows 23-26 top to bottom, with B25 set
# last.
$ws.Range("A21").Value = "ReturnNameExceptionalCase"
$ws.Range("B21").Value = "SER"
$ws.Range("C21").Value = "Any return which name ends with this word, must be skipped or ignored."

$ws.Range("A19").Value = "PathTemplateReviewSheet"
$ws.Range("B19").Value = "Data\Template_ReviewSheet.xlsx"

$ws.Range("A23").Value = "ReviewSheet_WorksheetName"
$ws.Range("B23").Value = "Template"

$ws.Range("A24").Value = "ReviewSheet_StartCellReturnsFailedTable"
$ws.Range("B24").Value = "G6"

$ws.Range("A25").Value = "ReviewSheet_RequiredColumns"

$ws.Range("A26").Value = "ReviewSheet_DateFormat"
$ws.Range("B26").Value = "MMMM yyyy"

$ws.Range("B25").Value = "Form Name,Legal Entity,Reason(s) Denied,Fixed?"

# Leave the sheet scrolled/selected where the author left it.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("B23").Select()
